$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-02 Tuesday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-04-03 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("211×9=1899", $true, $true, $false, $false, $false, $true, 1, $false, "624×8=4992", 2) | Out-Null
$d.Content.Find.Execute("800×7=5600", $true, $true, $false, $false, $false, $true, 1, $false, "161×9=1449", 2) | Out-Null
$d.Content.Find.Execute("622×5=3110", $true, $true, $false, $false, $false, $true, 1, $false, "628×7=4396", 2) | Out-Null
$d.Content.Find.Execute("665×3=1995", $true, $true, $false, $false, $false, $true, 1, $false, "341×8=2728", 2) | Out-Null
$d.Content.Find.Execute("796×5=3980", $true, $true, $false, $false, $false, $true, 1, $false, "980×9=8820", 2) | Out-Null
$d.Content.Find.Execute("121×3=363", $true, $true, $false, $false, $false, $true, 1, $false, "945×5=4725", 2) | Out-Null
$d.Content.Find.Execute("803×7=5621", $true, $true, $false, $false, $false, $true, 1, $false, "275×8=2200", 2) | Out-Null
$d.Content.Find.Execute("503×5=2515", $true, $true, $false, $false, $false, $true, 1, $false, "881×5=4405", 2) | Out-Null
$d.Content.Find.Execute("951×6=5706", $true, $true, $false, $false, $false, $true, 1, $false, "612×4=2448", 2) | Out-Null
$d.Content.Find.Execute("558×5=2790", $true, $true, $false, $false, $false, $true, 1, $false, "403×9=3627", 2) | Out-Null
$d.Content.Find.Execute("230×7=1610", $true, $true, $false, $false, $false, $true, 1, $false, "914×9=8226", 2) | Out-Null
$d.Content.Find.Execute("464×3=1392", $true, $true, $false, $false, $false, $true, 1, $false, "430×7=3010", 2) | Out-Null
$d.Content.Find.Execute("955×4=3820", $true, $true, $false, $false, $false, $true, 1, $false, "293×4=1172", 2) | Out-Null
$d.Content.Find.Execute("804×9=7236", $true, $true, $false, $false, $false, $true, 1, $false, "743×5=3715", 2) | Out-Null
$d.Content.Find.Execute("427×6=2562", $true, $true, $false, $false, $false, $true, 1, $false, "431×6=2586", 2) | Out-Null
$d.Content.Find.Execute("431×9=3879", $true, $true, $false, $false, $false, $true, 1, $false, "333×9=2997", 2) | Out-Null
$d.Content.Find.Execute("208×7=1456", $true, $true, $false, $false, $false, $true, 1, $false, "918×9=8262", 2) | Out-Null
$d.Content.Find.Execute("609×4=2436", $true, $true, $false, $false, $false, $true, 1, $false, "925×8=7400", 2) | Out-Null
$d.Content.Find.Execute("995×8=7960", $true, $true, $false, $false, $false, $true, 1, $false, "690×2=1380", 2) | Out-Null
$d.Content.Find.Execute("785×8=6280", $true, $true, $false, $false, $false, $true, 1, $false, "791×8=6328", 2) | Out-Null
$d.Content.Find.Execute("781×3=2343", $true, $true, $false, $false, $false, $true, 1, $false, "784×3=2352", 2) | Out-Null
$d.Content.Find.Execute("746×9=6714", $true, $true, $false, $false, $false, $true, 1, $false, "548×5=2740", 2) | Out-Null
$d.Content.Find.Execute("293×2=586", $true, $true, $false, $false, $false, $true, 1, $false, "699×4=2796", 2) | Out-Null
$d.Content.Find.Execute("767×7=5369", $true, $true, $false, $false, $false, $true, 1, $false, "193×3=579", 2) | Out-Null
$d.Content.Find.Execute("577×3=1731", $true, $true, $false, $false, $false, $true, 1, $false, "630×2=1260", 2) | Out-Null
